{"js": "// Insert a new paragraph, right after the paragraph that ends with\n// \"...So the average rate is 512.5 Mbps.\" and right before the paragraph\n// that begins with \"d) 4 points Now suppose ...\", containing the general\n// formula explanation that the author added to the exam answer key.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its distinctive text rather than by a\n// fixed index, so the script is resilient to minor structural differences.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"50 % of the requests can be delivered at 25 Mbps\") !== -1 &&\n      t.indexOf(\"512.5 Mbps\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the '...512.5 Mbps.' answer paragraph to anchor the insertion.\");\n}\n\nconst newParagraphText =\n  \"In general, suppose x of the requests can be delivered at 25 Mbps and \" +\n  \"1-x of the requests can be delivered at 1 Gbps. The average rate is \" +\n  \"x*25Mbps + (1-x)*1Gbps.\";\n\nanchor.insertParagraph(newParagraphText, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph, right after the paragraph that ends with\n# \"...So the average rate is 512.5 Mbps.\" and right before the paragraph\n# that begins with \"d) 4 points Now suppose ...\", containing the general\n# formula explanation that the author added to the exam answer key.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*50 % of the requests can be delivered at 25 Mbps*512.5 Mbps*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    Write-Output \"Could not find the '...512.5 Mbps.' answer paragraph to anchor the insertion.\"\n} else {\n    $newParagraphText = \"In general, suppose x of the requests can be delivered at 25 Mbps and 1-x of the requests can be delivered at 1 Gbps. The average rate is x*25Mbps + (1-x)*1Gbps.\"\n\n    $r = $target.Range\n    $r.InsertParagraphAfter()\n    $newParaRange = $d.Range($r.End, $r.End)\n    $newParaRange.InsertAfter($newParagraphText)\n}\n"}
